# Apply the "cryptos list" daily refresh (GitHub Actions scheduled update).
# Updates price (column D) and 1h volume/change % (column E) for most rows,
# and swaps the Chainlink/Uniswap rows (20 <-> 21) back to their new rank order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.816.68"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").Value = "3.845.38"
$ws.Range("E3").Value = "  +2.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.52"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.20"
$ws.Range("E6").Value = "  -2.99%  "

$ws.Range("D7").Value = "3.844.38"
$ws.Range("E7").Value = "  +2.40%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.47%  "

$ws.Range("E10").Value = "  -1.01%  "

$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.88"
$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("E14").Value = "  -2.02%  "

$ws.Range("D15").Value = "4.489.74"
$ws.Range("E15").Value = "  +2.45%  "

$ws.Range("D16").Value = "3.868.17"
$ws.Range("E16").Value = "  +3.31%  "

$ws.Range("D17").Value = "68.956.45"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.36"
$ws.Range("E20").Value = "  -0.12%  "

$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.12"
$ws.Range("E21").Value = "  -1.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "483.96"
$ws.Range("E22").Value = "  -1.97%  "

$ws.Range("E23").Value = "  -1.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000159"
$ws.Range("E24").Value = "  +3.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.93"
$ws.Range("E25").Value = "  -1.06%  "

$ws.Range("E26").Value = "  -2.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.09"
$ws.Range("E27").Value = "  -2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("E31").Value = "  -2.56%  "

$ws.Range("D32").Value = "3.992.42"
$ws.Range("E32").Value = "  +2.42%  "

$ws.Range("E33").Value = "  -4.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.19"
$ws.Range("E34").Value = "  +1.80%  "

$ws.Range("D35").Value = "3.792.67"
$ws.Range("E35").Value = "  +2.73%  "

$ws.Range("E36").Value = "  -1.60%  "

$ws.Range("E37").Value = "  +0.95%  "

$ws.Range("E38").Value = "  +2.96%  "

$ws.Range("E39").Value = "  -1.60%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  -1.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "437.62"
$ws.Range("E42").Value = "  +1.66%  "

$ws.Range("E43").Value = "  -1.81%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.48"
$ws.Range("E44").Value = "  -0.57%  "

$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("E47").Value = "  -1.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.38"
$ws.Range("E48").Value = "  +1.64%  "

$ws.Range("D49").Value = "2.834.45"
$ws.Range("E49").Value = "  +1.28%  "

$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.00"
$ws.Range("E51").Value = "  +11.36%  "
